$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches original inline-string cells)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.156.69'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.538.28'
$ws.Range("E3").Value = '  +2.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.49'
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.16'
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.535.80'
$ws.Range("E9").Value = '  +2.86%  '

$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.56'

$ws.Range("E13").Value = '  +0.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.953.87'
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.73'
$ws.Range("E15").Value = '  -1.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.083.38'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.522.85'
$ws.Range("E18").Value = '  +0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.22'
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.29'
$ws.Range("E20").Value = '  -1.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.44'
$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("E22").Value = '  +3.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.82'
$ws.Range("E23").Value = '  +1.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.90'
$ws.Range("E24").Value = '  +1.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.437'
$ws.Range("E25").Value = '  -4.88%  '

$ws.Range("E26").Value = '  +1.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.992'
$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.97'
$ws.Range("E28").Value = '  +3.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0779'
$ws.Range("E29").Value = '  +0.81%  '

$ws.Range("E30").Value = '  -0.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.22'
$ws.Range("E32").Value = '  -4.60%  '

$ws.Range("E33").Value = '  +8.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.81'
$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.68'
$ws.Range("E36").Value = '  +1.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.36'
$ws.Range("E37").Value = '  -2.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("E38").Value = '  -6.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.66'
$ws.Range("E39").Value = '  -4.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.96'
$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '297.18'
$ws.Range("E41").Value = '  -7.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.70'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.815'
$ws.Range("E43").Value = '  -2.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.603'
$ws.Range("E45").Value = '  +3.66%  '

$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0932'
$ws.Range("E47").Value = '  -1.13%  '

$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.33'

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0228'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0513'
$ws.Range("E51").Value = '  -2.42%  '
